$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-09 -> 2023-09-10) for every data row (rows 2 through 44).
for ($row = 2; $row -le 44; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value = 45179
    }
}
